$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the end time on row 59 (was 22:00 -> now 21:20)
$ws.Range("E59").Value = 0.88888888888888884

# 2. Insert 4 new blank rows before row 60, shifting the trailing
#    separator/summary rows down from 60-63 to 64-67.
$ws.Rows("60:63").Insert()

# 3. Fill in the new data rows 60-63 with the hand-over data.
$newRows = @(
    @{ Row = 60; Year = 2014; Month = 3; Day = 11; Start = 0.8125;              End = 0.96875 },
    @{ Row = 61; Year = 2014; Month = 3; Day = 12; Start = 0.35416666666666669; End = 0.41666666666666669 },
    @{ Row = 62; Year = 2014; Month = 3; Day = 12; Start = 0.48958333333333331; End = 0.52083333333333337 },
    @{ Row = 63; Year = 2014; Month = 3; Day = 12; Start = 0.63888888888888895; End = 0.75 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Year
    $ws.Cells.Item($row, 2).Value = $r.Month
    $ws.Cells.Item($row, 3).Value = $r.Day
    $ws.Cells.Item($row, 4).Value = $r.Start
    $ws.Cells.Item($row, 5).Value = $r.End
}

$ws.Range("F60:F63").FormulaR1C1 = "=(RC[-1]-RC[-2])*24*60"
$ws.Range("G60:G63").FormulaR1C1 = "=RC[-1]/60"

# 4. Move the selection onto the newly added block, matching the
#    author's cursor position after keying in the new rows.
$ws.Range("A64").Select() | Out-Null

$wb.Save()
